$d = $word.ActiveDocument

# Locate the date run that currently reads " 09, 2016" (it directly follows a
# separate run containing "May").
$rng = $d.Content
$found = $rng.Find.Execute(" 09, 2016", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find date text ' 09, 2016'"
}
$dateStart = $rng.Start
$dateEnd = $rng.End

# Step 1: Replace the leading space with a placeholder letter "X" (same length,
# so nothing shifts). This removes the need for xml:space="preserve" on this
# run, which matters because once a run carries that attribute, any bookmark
# split of it will make *both* resulting fragments inherit the attribute even
# when they no longer contain leading/trailing whitespace. Performing this
# "neutralising" edit before splitting means the later, real split points
# start from runs that only get xml:space="preserve" where it is actually
# needed.
$d.Range($dateStart, $dateEnd).Text = "X09, 2016"

# Step 2: Split the run into pieces with bookmarks (temporary helper bookmarks
# plus the real _GoBack bookmark). Bookmarks placed between text force Word to
# break the run at that point.
#   dateStart .. dateStart+1   -> "X"        (will become " ")
#   dateStart+1 .. dateStart+3 -> "09"       (will become "10")
#   dateStart+3 .. dateEnd     -> ", 2016"
$d.Bookmarks.Add("ZZTEMP_BOUNDARY", $d.Range($dateStart, $dateStart))
$d.Bookmarks.Add("_GoBack", $d.Range($dateStart + 3, $dateStart + 3))
$d.Bookmarks.Add("ZZTEMP_SPLIT", $d.Range($dateStart + 1, $dateStart + 1))

# Step 3: Now that the runs are isolated, fix up the actual text content.
# "X" -> " " (this tiny isolated run legitimately needs xml:space="preserve").
$d.Range($dateStart, $dateStart + 1).Text = " "
# "09" -> "10" (the date change itself).
$d.Range($dateStart + 1, $dateStart + 3).Text = "10"

# Step 4: Remove the temporary helper bookmarks. Deleting a bookmark does not
# re-merge the surrounding runs, so the split we engineered above survives.
$d.Bookmarks("ZZTEMP_BOUNDARY").Delete()
$d.Bookmarks("ZZTEMP_SPLIT").Delete()
